# Updated cryptos list on Wed Oct 25 11:08:54 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$value)
    # Force the cell to be stored as text so numeric-looking strings
    # (e.g. "224.65") are not silently converted into numbers, while
    # keeping the cell's original (default/General) number format.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "34.174.84"
$ws.Range("E2").Value = "  -0.75%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "1.780.60"
$ws.Range("E3").Value = "  -2.44%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.10%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "224.65"
$ws.Range("E5").Value = "  -2.29%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +0.20%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.06%  "

# Row 8 - Solana
Set-TextValue $ws.Range("D8") "31.93"
$ws.Range("E8").Value = "  +0.83%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -1.24%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -2.17%  "

# Row 11 - TRON
Set-TextValue $ws.Range("D11") "0.0931"
$ws.Range("E11").Value = "  +0.00%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D12") "2.034.67"
$ws.Range("E12").Value = "  -2.53%  "

# Row 13 - Chainlink
$ws.Range("E13").Value = "  +7.50%  "

# Row 14 - WrappedEther
Set-TextValue $ws.Range("D14") "1.775.37"
$ws.Range("E14").Value = "  -2.71%  "

# Row 15 - Polygon
Set-TextValue $ws.Range("D15") "0.627"
$ws.Range("E15").Value = "  -3.55%  "

# Row 16 - WrappedBTC
Set-TextValue $ws.Range("D16") "34.161.34"
$ws.Range("E16").Value = "  -0.55%  "

# Row 17 - Polkadot
Set-TextValue $ws.Range("D17") "4.21"
$ws.Range("E17").Value = "  -1.89%  "

# Row 18 - Litecoin
$ws.Range("E18").Value = "  -1.40%  "

# Row 19 - BitcoinCash
Set-TextValue $ws.Range("D19") "255.00"
$ws.Range("E19").Value = "  -1.51%  "

# Row 20 - ShibaInu
$ws.Range("E20").Value = "  -1.93%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  +0.13%  "

# Row 22 - Avalanche
Set-TextValue $ws.Range("D22") "10.36"

# Row 23 - Uniswap
Set-TextValue $ws.Range("D23") "4.19"
$ws.Range("E23").Value = "  -3.86%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  -4.21%  "

# Row 25 - Monero
Set-TextValue $ws.Range("D25") "157.28"
$ws.Range("E25").Value = "  -1.04%  "

# Row 26 - EthereumClassic
Set-TextValue $ws.Range("D26") "16.39"
$ws.Range("E26").Value = "  -1.56%  "

# Row 27 - Cosmos
$ws.Range("E27").Value = "  -1.81%  "

# Row 28 - Stellar
$ws.Range("E28").Value = "  -1.39%  "

# Row 29 - BinanceUSD
Set-TextValue $ws.Range("D29") "1.00"
$ws.Range("E29").Value = "  +0.13%  "

# Row 30 - Filecoin
Set-TextValue $ws.Range("D30") "3.77"
$ws.Range("E30").Value = "  -3.22%  "

# Row 31 - Hedera
Set-TextValue $ws.Range("D31") "0.0513"
$ws.Range("E31").Value = "  -1.69%  "

# Row 32 - PancakeSwap
$ws.Range("E32").Value = "  -1.89%  "

# Row 33 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D33") "3.58"
$ws.Range("E33").Value = "  +0.18%  "

# Row 34 - LidoDAOToken
Set-TextValue $ws.Range("D34") "1.86"
$ws.Range("E34").Value = "  +4.18%  "

# Row 35 - Maker
Set-TextValue $ws.Range("D35") "1.438.64"
$ws.Range("E35").Value = "  -7.40%  "

# Row 36 - TrustWalletToken
$ws.Range("E36").Value = "  -3.22%  "

# Row 37 - VeChain
$ws.Range("E37").Value = "  -1.61%  "

# Row 38 - ImmutableX
Set-TextValue $ws.Range("D38") "0.623"
$ws.Range("E38").Value = "  -2.26%  "

# Row 39 - was MXToken, now Aave
$ws.Range("B39").Value = "Aave"
$ws.Range("C39").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D39") "82.82"
$ws.Range("E39").Value = "  -2.19%  "

# Row 40 - was Aave, now MXToken
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D40") "2.84"
$ws.Range("E40").Value = "  +0.59%  "

# Row 41 - HuobiToken
$ws.Range("E41").Value = "  +0.81%  "

# Row 42 - ARBITRUM
Set-TextValue $ws.Range("D42") "0.889"
$ws.Range("E42").Value = "  -3.24%  "

# Row 43 - RenderToken
$ws.Range("E43").Value = "  -5.32%  "

# Row 44 - Kaspa
$ws.Range("E44").Value = "  -2.68%  "

# Row 45 - WEMIXToken
$ws.Range("E45").Value = "  -1.97%  "

# Row 46 - FraxShare
$ws.Range("E46").Value = "  +0.47%  "

# Row 47 - RocketPoolETH
Set-TextValue $ws.Range("D47") "1.936.43"
$ws.Range("E47").Value = "  -2.74%  "

# Row 48 - InjectiveProtocol
Set-TextValue $ws.Range("D48") "12.20"
$ws.Range("E48").Value = "  -2.02%  "

# Row 49 - PaxDollar
$ws.Range("E49").Value = "  +0.07%  "

# Row 50 - Quant
Set-TextValue $ws.Range("D50") "98.46"
$ws.Range("E50").Value = "  +0.75%  "

# Row 51 - BitcoinSV
Set-TextValue $ws.Range("D51") "49.47"
$ws.Range("E51").Value = "  -6.48%  "
